$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slightly taller rows for the existing data rows (2-10), matching the
# re-saved workbook's row metrics.
$ws.Range("A2:B10").RowHeight = 18.75

# New "N/A" label row appended below the existing walk-group table.
$ws.Range("A12").Value = "N/A"
$ws.Range("B12").Value = "N/A"

# Match the formatting of the row immediately above (same body style/height).
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12:B12").RowHeight = 18
